# Scheduled-runner refresh of Leve market/profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the per-job Leve tables. Values only; no
# formulas or formatting are touched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 974.0952
$ws.Range("I19").Value = 596.8333
$ws.Range("J19").Value = 1125
$ws.Range("K19").Value = 596.8333
$ws.Range("L19").Value = 1125
$ws.Range("M19").Value = -421.8333
$ws.Range("N19").Value = -1475
$ws.Range("H86").Value = 7627.4546
$ws.Range("I86").Value = 2667.6667
$ws.Range("J86").Value = 9487.375
$ws.Range("K86").Value = 2667.6667
$ws.Range("L86").Value = 9487.375
$ws.Range("M86").Value = -1544.6667
$ws.Range("N86").Value = -11733.375
$ws.Range("H89").Value = 7627.4546
$ws.Range("I89").Value = 2667.6667
$ws.Range("J89").Value = 9487.375
$ws.Range("K89").Value = 13338.3335
$ws.Range("L89").Value = 47436.875
$ws.Range("M89").Value = -7722.333500000001
$ws.Range("N89").Value = -58668.875
$ws.Range("H94").Value = 83334456
$ws.Range("I94").Value = 1228.4546
$ws.Range("K94").Value = 1228.4546
$ws.Range("M94").Value = -777.4546
$ws.Range("H103").Value = 1225.1333
$ws.Range("I103").Value = 1562
$ws.Range("J103").Value = 719.8333
$ws.Range("K103").Value = 4686
$ws.Range("L103").Value = 2159.4999
$ws.Range("M103").Value = -4100
$ws.Range("N103").Value = -3331.4999
$ws.Range("H107").Value = 376.78262
$ws.Range("I107").Value = 333
$ws.Range("K107").Value = 333
$ws.Range("M107").Value = 1587
$ws.Range("H127").Value = 1816.1578
$ws.Range("J127").Value = 2069.4688
$ws.Range("L127").Value = 6208.4064
$ws.Range("N127").Value = -16128.4064
$ws.Range("H129").Value = 2712.07
$ws.Range("J129").Value = 1017.7959
$ws.Range("L129").Value = 3053.3877
$ws.Range("N129").Value = -13053.3877
$ws.Range("H132").Value = 2843820
$ws.Range("I132").Value = 2944087.5
$ws.Range("K132").Value = 8832262.5
$ws.Range("M132").Value = -8829732.5
$ws.Range("H137").Value = 2117.68
$ws.Range("I137").Value = 1402.3889
$ws.Range("J137").Value = 3957
$ws.Range("K137").Value = 4207.1667
$ws.Range("L137").Value = 11871
$ws.Range("M137").Value = -1657.1667
$ws.Range("N137").Value = -16971
$ws.Range("H138").Value = 8505
$ws.Range("I138").Value = 1488.1072
$ws.Range("J138").Value = 47799.6
$ws.Range("K138").Value = 4464.321599999999
$ws.Range("L138").Value = 143398.8
$ws.Range("M138").Value = 675.6784000000007
$ws.Range("N138").Value = -153678.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25538.535
$ws.Range("I32").Value = 4224.721
$ws.Range("K32").Value = 4224.721
$ws.Range("M32").Value = -3937.721
$ws.Range("H61").Value = 1543.6897
$ws.Range("I61").Value = 1024.5555
$ws.Range("J61").Value = 2393.182
$ws.Range("K61").Value = 1024.5555
$ws.Range("L61").Value = 2393.182
$ws.Range("M61").Value = -812.5554999999999
$ws.Range("N61").Value = -2817.182
$ws.Range("H110").Value = 143157840
$ws.Range("I110").Value = 167017150
$ws.Range("J110").Value = 1999
$ws.Range("K110").Value = 167017150
$ws.Range("L110").Value = 1999
$ws.Range("M110").Value = -167015105
$ws.Range("N110").Value = -6089
$ws.Range("H136").Value = 1543.6897
$ws.Range("I136").Value = 1024.5555
$ws.Range("J136").Value = 2393.182
$ws.Range("K136").Value = 3073.6665
$ws.Range("L136").Value = 7179.545999999999
$ws.Range("M136").Value = -523.6664999999998
$ws.Range("N136").Value = -12279.546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2927.1052
$ws.Range("I134").Value = 2657.25
$ws.Range("K134").Value = 7971.75
$ws.Range("M134").Value = -5436.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 20000
$ws.Range("I32").Value = 20000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 20000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -19684
$ws.Range("N32").ClearContents()
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1595.3429
$ws.Range("I5").Value = 768.5238000000001
$ws.Range("J5").Value = 2835.5715
$ws.Range("K5").Value = 2305.5714
$ws.Range("L5").Value = 8506.7145
$ws.Range("M5").Value = -2193.5714
$ws.Range("N5").Value = -8730.7145
$ws.Range("H113").Value = 633.9231
$ws.Range("I113").Value = 578.26666
$ws.Range("J113").Value = 668.7083
$ws.Range("K113").Value = 1734.79998
$ws.Range("L113").Value = 2006.1249
$ws.Range("M113").Value = 435.20002
$ws.Range("N113").Value = -6346.1249
$ws.Range("H131").Value = 848.27
$ws.Range("J131").Value = 873.7578999999999
$ws.Range("L131").Value = 2621.2737
$ws.Range("N131").Value = -12701.2737
$ws.Range("H132").Value = 2313.9285
$ws.Range("J132").Value = 3110.5557
$ws.Range("L132").Value = 27995.0013
$ws.Range("N132").Value = -33055.0013
$ws.Range("H135").Value = 1595.3429
$ws.Range("I135").Value = 768.5238000000001
$ws.Range("J135").Value = 2835.5715
$ws.Range("K135").Value = 6916.7142
$ws.Range("L135").Value = 25520.1435
$ws.Range("M135").Value = -4381.7142
$ws.Range("N135").Value = -30590.1435
$ws.Range("H138").Value = 8312.0625
$ws.Range("I138").Value = 9971.666999999999
$ws.Range("K138").Value = 29915.001
$ws.Range("M138").Value = -24775.001
$ws.Range("H141").Value = 2251.125
$ws.Range("I141").Value = 1715.5714
$ws.Range("J141").Value = 6000
$ws.Range("K141").Value = 5146.7142
$ws.Range("L141").Value = 18000
$ws.Range("M141").Value = 33.28579999999965
$ws.Range("N141").Value = -28360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 70004
$ws.Range("J19").Value = 70004
$ws.Range("L19").Value = 70004
$ws.Range("N19").Value = -70580
$ws.Range("H96").Value = 59800
$ws.Range("J96").Value = 59800
$ws.Range("L96").Value = 59800
$ws.Range("N96").Value = -65292
$ws.Range("H132").Value = 3616.52
$ws.Range("I132").Value = 2300.8462
$ws.Range("K132").Value = 6902.5386
$ws.Range("M132").Value = -4372.5386

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 31595
$ws.Range("J6").Value = 31595
$ws.Range("L6").Value = 31595
$ws.Range("N6").Value = -31819
$ws.Range("H16").Value = 50893.05
$ws.Range("I16").Value = 59404.06
$ws.Range("J16").Value = 2664
$ws.Range("K16").Value = 59404.06
$ws.Range("L16").Value = 2664
$ws.Range("M16").Value = -59234.06
$ws.Range("N16").Value = -3004
$ws.Range("H95").Value = 25000
$ws.Range("J95").Value = 25000
$ws.Range("L95").Value = 25000
$ws.Range("N95").Value = -30492
$ws.Range("H96").Value = 19500
$ws.Range("J96").Value = 19500
$ws.Range("L96").Value = 19500
$ws.Range("N96").Value = -24992

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 77710.46000000001
$ws.Range("I107").Value = 550.6667
$ws.Range("J107").Value = 100858.4
$ws.Range("K107").Value = 1652.0001
$ws.Range("L107").Value = 302575.2
$ws.Range("M107").Value = 267.9999
$ws.Range("N107").Value = -306415.2
$ws.Range("H113").Value = 713.7778
$ws.Range("I113").Value = 404
$ws.Range("J113").Value = 1333.3334
$ws.Range("K113").Value = 1212
$ws.Range("L113").Value = 4000.0002
$ws.Range("M113").Value = 958
$ws.Range("N113").Value = -8340.0002
$ws.Range("H124").Value = 40330
$ws.Range("J124").Value = 40330
$ws.Range("L124").Value = 40330
$ws.Range("N124").Value = -50150
$ws.Range("H133").Value = 63310
$ws.Range("J133").Value = 63310
$ws.Range("L133").Value = 63310
$ws.Range("N133").Value = -73430
